$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data released 2020-11-09: updated nombre_aides (col C) and montant_total (col D)
# values for a set of region/categorie_juridique rows. Values are kept as text
# (matching the workbook's existing inline-string storage) to avoid floating
# point rounding / trailing-zero loss, and the original cell style is restored
# after the write so no formatting changes are introduced.
$data = @(
    @{Row=2; C="1707"; D="4286914.49"},
    @{Row=4; C="1282"; D="6764472.71"},
    @{Row=6; C="882"; D="4128790.54"},
    @{Row=9; C="265"; D="905212.47"},
    @{Row=11; C="456"; D="2744003.28"},
    @{Row=23; C="439"; D="2689998.61"},
    @{Row=31; C="438"; D="1290569.11"},
    @{Row=33; C="816"; D="5340051.92"},
    @{Row=35; C="545"; D="2953342.32"},
    @{Row=38; C="585"; D="1536372.32"},
    @{Row=39; C="284"; D="1293780.04"},
    @{Row=40; C="273"; D="924520.72"},
    @{Row=43; C="481"; D="1758556.94"},
    @{Row=44; C="191"; D="1158797.39"},
    @{Row=45; C="272"; D="1328145.19"},
    @{Row=47; C="20"; D="122181.23"},
    @{Row=48; C="762"; D="2278051.35"},
    @{Row=50; C="995"; D="6364799.81"},
    @{Row=51; C="722"; D="3980122.37"},
    @{Row=54; C="10489"; D="31159641.50"},
    @{Row=57; C="55"; D="424009.00"},
    @{Row=58; C="6948"; D="35542073.29"},
    @{Row=59; C="23"; D="253000.00"},
    @{Row=60; C="6810"; D="29468918.34"},
    @{Row=62; C="136"; D="686657.46"},
    @{Row=79; C="444"; D="1322392.80"},
    @{Row=81; C="1181"; D="7720488.19"},
    @{Row=82; C="634"; D="3637090.36"},
    @{Row=83; C="48"; D="185078.00"},
    @{Row=85; C="931"; D="2554538.79"},
    @{Row=88; C="1329"; D="7953593.54"},
    @{Row=90; C="942"; D="5112684.28"},
    @{Row=91; C="35"; D="111206.00"},
    @{Row=93; C="255"; D="648350.00"},
    @{Row=100; C="1374"; D="3479754.28"},
    @{Row=103; C="1551"; D="7643841.24"},
    @{Row=105; C="1495"; D="6670911.06"}
)

foreach ($item in $data) {
    $r = $item.Row

    $cCell = $ws.Cells.Item($r, 3)
    $cStyle = $cCell.Style
    $cCell.NumberFormat = "@"
    $cCell.Value = $item.C
    $cCell.Style = $cStyle

    $dCell = $ws.Cells.Item($r, 4)
    $dStyle = $dCell.Style
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.D
    $dCell.Style = $dStyle
}
